# Links.xlsx update:
#  - swap the 9 product-photo links (column A) and their numeric ids
#    (column B) for a new batch pulled from a different supplier site
#  - add a helper column C that mirrors column A's "nome" text so a
#    later lookup/formula has a local copy next to the id
#  - the sheet now uses A1:C.. instead of A1:B.. and is set up for
#    printing on A4 portrait paper

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New links (column A / shared with column C) and their matching ids
# (column B), row-for-row in the same order as the old data they replace.
$links = @(
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=21989",
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=50261",
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=650158",
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=50193",
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=50282",
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=60124",
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=21503",
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=50127",
    "https://www.ideia2001.com.br/catmobile/FotoMobRetArq.asp?cerq=226&n=50195"
)

$ids = @(21989, 50261, 650158, 50193, 50282, 60124, 21503, 50127, 50195)

for ($i = 0; $i -lt $links.Length; $i++) {
    $row = $i + 2   # data starts at row 2 (row 1 is the "link"/"nome" header)
    $ws.Cells.Item($row, 1).Value = $links[$i]   # A: link
    $ws.Cells.Item($row, 2).Value = $ids[$i]     # B: nome (numeric id)
    $ws.Cells.Item($row, 3).Value = $links[$i]   # C: new helper column, mirrors A
}

# Page setup: print on A4, portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
